# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Sheet "展览" (sheet1 / index 1)
$ws1 = $excel.ActiveWorkbook.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1461
$ws1.Range("F4").Value = 1748
$ws1.Range("F6").Value = 142
$ws1.Range("F8").Value = 33
$ws1.Range("F10").Value = 551
$ws1.Range("F13").Value = 145
$ws1.Range("F18").Value = 4636
$ws1.Range("F20").Value = 817
$ws1.Range("F21").Value = 99
$ws1.Range("F22").Value = 2185
$ws1.Range("F24").Value = 14
$ws1.Range("F25").Value = 2050

# Sheet "演出" (sheet2)
$ws2 = $excel.ActiveWorkbook.Worksheets.Item("演出")
$ws2.Range("F2").Value = 74

# Sheet "全部类型" (sheet4)
$ws4 = $excel.ActiveWorkbook.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1461
$ws4.Range("F4").Value = 1748
$ws4.Range("F6").Value = 142
$ws4.Range("F8").Value = 33
$ws4.Range("F10").Value = 551
$ws4.Range("F13").Value = 145
$ws4.Range("F18").Value = 4636
$ws4.Range("F19").Value = 74
$ws4.Range("F22").Value = 817
$ws4.Range("F23").Value = 99
$ws4.Range("F24").Value = 2185
$ws4.Range("F26").Value = 14
$ws4.Range("F27").Value = 2050
